$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.578.52"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.424.13"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'306.51"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").Value = "'97.06"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "'35.09"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "'18.48"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "2.795.03"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "2.416.14"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "'0.827"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").Value = "43.608.39"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "'12.04"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'237.73"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "'2.26"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'24.98"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "'9.43"
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("D30").Value = "'32.31"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  +18.27%  "
$ws.Range("D32").Value = "'18.46"
$ws.Range("E32").Value = "  +6.85%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'0.0750"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.90"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'129.97"
$ws.Range("E37").Value = "  +27.82%  "
$ws.Range("D38").Value = "'2.93"
$ws.Range("E38").Value = "  +5.04%  "
$ws.Range("D39").Value = "'4.38"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "'2.27"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "'21.12"
$ws.Range("E42").Value = "  -6.97%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "'0.0283"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "2.657.59"
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").Value = "'52.62"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").Value = "'72.39"
$ws.Range("E51").Value = "  +0.06%  "

$forceCells = @("D4","D5","D6","D10","D13","D17","D19","D20","D23","D24","D27","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D42","D44","D50","D51")
foreach ($c in $forceCells) {
    $ws.Range($c).Style = "Normal"
}

Write-Host "Applied cryptos update"